$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "1.000", "4.260") are preserved exactly as text, matching the
# original workbook where every data cell is stored as an inline string.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.447.37"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.871.24"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "0.7075"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").Value = "243.73"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "0.3164"
$ws.Range("E8").Value = "  +0.79%  "

$ws.Range("D9").Value = "0.07885"
$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("D10").Value = "24.59"
$ws.Range("E10").Value = "  -2.03%  "

$ws.Range("D11").Value = "0.07994"
$ws.Range("E11").Value = "  -4.02%  "

$ws.Range("D12").Value = "1.882.76"
$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("D13").Value = "5.224"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "0.7052"
$ws.Range("E15").Value = "  -1.77%  "

$ws.Range("D16").Value = "6.516"
$ws.Range("E16").Value = "  +2.54%  "

$ws.Range("D17").Value = "29.472.72"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").Value = "0.000008364"
$ws.Range("E18").Value = "  -3.70%  "

$ws.Range("D19").Value = "257.64"
$ws.Range("E19").Value = "  +6.05%  "

$ws.Range("D20").Value = "2.124.06"
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").Value = "13.23"
$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").Value = "7.639"
$ws.Range("E23").Value = "  -2.88%  "

$ws.Range("D24").Value = "0.9999"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "0.1563"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").Value = "9.086"
$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").Value = "160.88"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("D28").Value = "18.94"
$ws.Range("E28").Value = "  +1.90%  "

$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").Value = "4.344"
$ws.Range("E30").Value = "  -2.07%  "

$ws.Range("D31").Value = "4.260"
$ws.Range("E31").Value = "  -2.41%  "

$ws.Range("D32").Value = "1.211"
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").Value = "0.05320"
$ws.Range("E33").Value = "  -1.46%  "

$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("D35").Value = "1.176"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "0.7485"
$ws.Range("E36").Value = "  -3.50%  "

$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  +0.97%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "1.266.99"
$ws.Range("E39").Value = "  -0.35%  "

$ws.Range("D40").Value = "2.757"
$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("D41").Value = "0.9049"
$ws.Range("E41").Value = "  -1.59%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.992"
$ws.Range("E42").Value = "  -8.51%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "71.85"
$ws.Range("E43").Value = "  -3.65%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "108.04"
$ws.Range("E44").Value = "  -5.05%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000130"
$ws.Range("E45").Value = "  +2.30%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9991"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").Value = "2.026.60"
$ws.Range("E47").Value = "  -0.59%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.795"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.5194"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").Value = "9.547"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").Value = "0.4332"
$ws.Range("E51").Value = "  -0.99%  "
